$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBCS3A3DV%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T200034Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDjxIHNpE7YwhIVzSwINqU1JKJWBKmFHslxXG%2B1HYOVNwIhALt4e1jhgjcn6fgzFDyDu3w4pkKamJ4dYU6KRFVuWtn8Kt0DCBwQBBoMNzgxMjA0OTQyMjQ0Igw8K0FKPD5wdCqMNo0qugPV2k2p8kv6vVW49wi%2FVk1v8utmS2HoMcRlMGAbytQ7LqFro03sqDoCVRZwb5It31S2RIJccfRZQEDNhr0k9QFpUtIbPJCrJ%2FYJyH1YKzLLIQc06SfrfrzqTbp06ZA6vjxPcY8dk7urOIH5ns2y1g1JHgE2egkFbFt%2BZ0sHEg8QWTLkSQk2kfHVjaszfQJ5DHB9tSBD6LAmNiGHJxEwCP%2FFa7cyieDLv5%2BtwfirJFH%2Fedb8d%2FJhTyrF7%2FizaLrOAP1Ek7Ba3LM9p0tLbB%2BQicvRpnXQgmwWUP4zGIWgCpZq7aahP405Ht%2BSYLUYxdVPs%2FdpRgB25KH7MTJflrPgo8kTwGK5780tyl5DcLYQ6oLi54mcrQymRVKADYrHl%2FUdtY6KXG9vM%2BA3meZi4W%2B1wJLXqKkrSb2%2FjJnaB1A77Om5StksZpMNmXCMWQ%2BLhsFvA%2BWcCo1DJGDRN5b1Q4mcBDQ9p6Y8Be5H8G03VOGsaz0M0J6BlEWVWoqSybYEipbrlos7%2FZvAVTmfzRAkfZcRQhuW%2FjCJYQ5xIg%2ByEtuTabai%2B0%2FkLEUhCeg5OkPEnjZRUxtsc2BeO3CPHCUgMKnIwL8GOqQBPR2ZPF%2BAn2N%2FAJ6xflLVSgZeKaDZERMT5Rq5amawGKejbyICI0CYtevbPzg93Sf9%2FcL4jrUO5NzPf0UbTZn6kHFezvSrlRUWZnr6NAMMaqMT4SsxaSA6%2BZo598UwTaOIjGA4fPti19zYRr4YMhi1Rth4ew0vAjChRwN4bCjj2ziiCxyOHEgU2eDosny4bMPL8eABhd8KitY73d5TNeinH3gtDgM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2004%3A00%3A34%20GMT&X-Amz-Signature=bd9711177cc8a6532563f98ebdbbf193ceb8d4c781f10ecfa33565834f0dee33'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBCS3A3DV%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T200034Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDjxIHNpE7YwhIVzSwINqU1JKJWBKmFHslxXG%2B1HYOVNwIhALt4e1jhgjcn6fgzFDyDu3w4pkKamJ4dYU6KRFVuWtn8Kt0DCBwQBBoMNzgxMjA0OTQyMjQ0Igw8K0FKPD5wdCqMNo0qugPV2k2p8kv6vVW49wi%2FVk1v8utmS2HoMcRlMGAbytQ7LqFro03sqDoCVRZwb5It31S2RIJccfRZQEDNhr0k9QFpUtIbPJCrJ%2FYJyH1YKzLLIQc06SfrfrzqTbp06ZA6vjxPcY8dk7urOIH5ns2y1g1JHgE2egkFbFt%2BZ0sHEg8QWTLkSQk2kfHVjaszfQJ5DHB9tSBD6LAmNiGHJxEwCP%2FFa7cyieDLv5%2BtwfirJFH%2Fedb8d%2FJhTyrF7%2FizaLrOAP1Ek7Ba3LM9p0tLbB%2BQicvRpnXQgmwWUP4zGIWgCpZq7aahP405Ht%2BSYLUYxdVPs%2FdpRgB25KH7MTJflrPgo8kTwGK5780tyl5DcLYQ6oLi54mcrQymRVKADYrHl%2FUdtY6KXG9vM%2BA3meZi4W%2B1wJLXqKkrSb2%2FjJnaB1A77Om5StksZpMNmXCMWQ%2BLhsFvA%2BWcCo1DJGDRN5b1Q4mcBDQ9p6Y8Be5H8G03VOGsaz0M0J6BlEWVWoqSybYEipbrlos7%2FZvAVTmfzRAkfZcRQhuW%2FjCJYQ5xIg%2ByEtuTabai%2B0%2FkLEUhCeg5OkPEnjZRUxtsc2BeO3CPHCUgMKnIwL8GOqQBPR2ZPF%2BAn2N%2FAJ6xflLVSgZeKaDZERMT5Rq5amawGKejbyICI0CYtevbPzg93Sf9%2FcL4jrUO5NzPf0UbTZn6kHFezvSrlRUWZnr6NAMMaqMT4SsxaSA6%2BZo598UwTaOIjGA4fPti19zYRr4YMhi1Rth4ew0vAjChRwN4bCjj2ziiCxyOHEgU2eDosny4bMPL8eABhd8KitY73d5TNeinH3gtDgM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2004%3A00%3A34%20GMT&X-Amz-Signature=2297b65a605e20a726a86b617874492f23ddb3fbbf4a87b5cbaea51776f73f6c'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBCS3A3DV%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T200034Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDjxIHNpE7YwhIVzSwINqU1JKJWBKmFHslxXG%2B1HYOVNwIhALt4e1jhgjcn6fgzFDyDu3w4pkKamJ4dYU6KRFVuWtn8Kt0DCBwQBBoMNzgxMjA0OTQyMjQ0Igw8K0FKPD5wdCqMNo0qugPV2k2p8kv6vVW49wi%2FVk1v8utmS2HoMcRlMGAbytQ7LqFro03sqDoCVRZwb5It31S2RIJccfRZQEDNhr0k9QFpUtIbPJCrJ%2FYJyH1YKzLLIQc06SfrfrzqTbp06ZA6vjxPcY8dk7urOIH5ns2y1g1JHgE2egkFbFt%2BZ0sHEg8QWTLkSQk2kfHVjaszfQJ5DHB9tSBD6LAmNiGHJxEwCP%2FFa7cyieDLv5%2BtwfirJFH%2Fedb8d%2FJhTyrF7%2FizaLrOAP1Ek7Ba3LM9p0tLbB%2BQicvRpnXQgmwWUP4zGIWgCpZq7aahP405Ht%2BSYLUYxdVPs%2FdpRgB25KH7MTJflrPgo8kTwGK5780tyl5DcLYQ6oLi54mcrQymRVKADYrHl%2FUdtY6KXG9vM%2BA3meZi4W%2B1wJLXqKkrSb2%2FjJnaB1A77Om5StksZpMNmXCMWQ%2BLhsFvA%2BWcCo1DJGDRN5b1Q4mcBDQ9p6Y8Be5H8G03VOGsaz0M0J6BlEWVWoqSybYEipbrlos7%2FZvAVTmfzRAkfZcRQhuW%2FjCJYQ5xIg%2ByEtuTabai%2B0%2FkLEUhCeg5OkPEnjZRUxtsc2BeO3CPHCUgMKnIwL8GOqQBPR2ZPF%2BAn2N%2FAJ6xflLVSgZeKaDZERMT5Rq5amawGKejbyICI0CYtevbPzg93Sf9%2FcL4jrUO5NzPf0UbTZn6kHFezvSrlRUWZnr6NAMMaqMT4SsxaSA6%2BZo598UwTaOIjGA4fPti19zYRr4YMhi1Rth4ew0vAjChRwN4bCjj2ziiCxyOHEgU2eDosny4bMPL8eABhd8KitY73d5TNeinH3gtDgM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2004%3A00%3A34%20GMT&X-Amz-Signature=c9e224b0ee94f33c6d1c134c2990491b21fc704b6fc6cbb7905544606a95f323'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBCS3A3DV%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T200034Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDjxIHNpE7YwhIVzSwINqU1JKJWBKmFHslxXG%2B1HYOVNwIhALt4e1jhgjcn6fgzFDyDu3w4pkKamJ4dYU6KRFVuWtn8Kt0DCBwQBBoMNzgxMjA0OTQyMjQ0Igw8K0FKPD5wdCqMNo0qugPV2k2p8kv6vVW49wi%2FVk1v8utmS2HoMcRlMGAbytQ7LqFro03sqDoCVRZwb5It31S2RIJccfRZQEDNhr0k9QFpUtIbPJCrJ%2FYJyH1YKzLLIQc06SfrfrzqTbp06ZA6vjxPcY8dk7urOIH5ns2y1g1JHgE2egkFbFt%2BZ0sHEg8QWTLkSQk2kfHVjaszfQJ5DHB9tSBD6LAmNiGHJxEwCP%2FFa7cyieDLv5%2BtwfirJFH%2Fedb8d%2FJhTyrF7%2FizaLrOAP1Ek7Ba3LM9p0tLbB%2BQicvRpnXQgmwWUP4zGIWgCpZq7aahP405Ht%2BSYLUYxdVPs%2FdpRgB25KH7MTJflrPgo8kTwGK5780tyl5DcLYQ6oLi54mcrQymRVKADYrHl%2FUdtY6KXG9vM%2BA3meZi4W%2B1wJLXqKkrSb2%2FjJnaB1A77Om5StksZpMNmXCMWQ%2BLhsFvA%2BWcCo1DJGDRN5b1Q4mcBDQ9p6Y8Be5H8G03VOGsaz0M0J6BlEWVWoqSybYEipbrlos7%2FZvAVTmfzRAkfZcRQhuW%2FjCJYQ5xIg%2ByEtuTabai%2B0%2FkLEUhCeg5OkPEnjZRUxtsc2BeO3CPHCUgMKnIwL8GOqQBPR2ZPF%2BAn2N%2FAJ6xflLVSgZeKaDZERMT5Rq5amawGKejbyICI0CYtevbPzg93Sf9%2FcL4jrUO5NzPf0UbTZn6kHFezvSrlRUWZnr6NAMMaqMT4SsxaSA6%2BZo598UwTaOIjGA4fPti19zYRr4YMhi1Rth4ew0vAjChRwN4bCjj2ziiCxyOHEgU2eDosny4bMPL8eABhd8KitY73d5TNeinH3gtDgM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2004%3A00%3A34%20GMT&X-Amz-Signature=f42cbaadabbd3043cffa93424ed47a336414d1f3b64ff5e05ce3d3cd37236ff8'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBCS3A3DV%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T200035Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDjxIHNpE7YwhIVzSwINqU1JKJWBKmFHslxXG%2B1HYOVNwIhALt4e1jhgjcn6fgzFDyDu3w4pkKamJ4dYU6KRFVuWtn8Kt0DCBwQBBoMNzgxMjA0OTQyMjQ0Igw8K0FKPD5wdCqMNo0qugPV2k2p8kv6vVW49wi%2FVk1v8utmS2HoMcRlMGAbytQ7LqFro03sqDoCVRZwb5It31S2RIJccfRZQEDNhr0k9QFpUtIbPJCrJ%2FYJyH1YKzLLIQc06SfrfrzqTbp06ZA6vjxPcY8dk7urOIH5ns2y1g1JHgE2egkFbFt%2BZ0sHEg8QWTLkSQk2kfHVjaszfQJ5DHB9tSBD6LAmNiGHJxEwCP%2FFa7cyieDLv5%2BtwfirJFH%2Fedb8d%2FJhTyrF7%2FizaLrOAP1Ek7Ba3LM9p0tLbB%2BQicvRpnXQgmwWUP4zGIWgCpZq7aahP405Ht%2BSYLUYxdVPs%2FdpRgB25KH7MTJflrPgo8kTwGK5780tyl5DcLYQ6oLi54mcrQymRVKADYrHl%2FUdtY6KXG9vM%2BA3meZi4W%2B1wJLXqKkrSb2%2FjJnaB1A77Om5StksZpMNmXCMWQ%2BLhsFvA%2BWcCo1DJGDRN5b1Q4mcBDQ9p6Y8Be5H8G03VOGsaz0M0J6BlEWVWoqSybYEipbrlos7%2FZvAVTmfzRAkfZcRQhuW%2FjCJYQ5xIg%2ByEtuTabai%2B0%2FkLEUhCeg5OkPEnjZRUxtsc2BeO3CPHCUgMKnIwL8GOqQBPR2ZPF%2BAn2N%2FAJ6xflLVSgZeKaDZERMT5Rq5amawGKejbyICI0CYtevbPzg93Sf9%2FcL4jrUO5NzPf0UbTZn6kHFezvSrlRUWZnr6NAMMaqMT4SsxaSA6%2BZo598UwTaOIjGA4fPti19zYRr4YMhi1Rth4ew0vAjChRwN4bCjj2ziiCxyOHEgU2eDosny4bMPL8eABhd8KitY73d5TNeinH3gtDgM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2004%3A00%3A35%20GMT&X-Amz-Signature=3aa4a1ead82b9a92587e2daf91cf14f72c105ade4ddbe40fa15180bbaca75b0e'
